$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2164.0598
$ws.Range("I137").Value = 1801.9149
$ws.Range("J137").Value = 3015.1
$ws.Range("K137").Value = 5405.7447
$ws.Range("L137").Value = 9045.299999999999
$ws.Range("M137").Value = -2855.7447
$ws.Range("N137").Value = -14145.3
$ws.Range("H138").Value = 2287.9824
$ws.Range("I138").Value = 762.9666999999999
$ws.Range("J138").Value = 3982.4443
$ws.Range("K138").Value = 2288.9001
$ws.Range("L138").Value = 11947.3329
$ws.Range("M138").Value = 2851.0999
$ws.Range("N138").Value = -22227.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()
$ws.Range("H30").Value = 1503
$ws.Range("I30").Value = 1009
$ws.Range("J30").Value = 1750
$ws.Range("K30").Value = 1009
$ws.Range("L30").Value = 1750
$ws.Range("M30").Value = -859
$ws.Range("N30").Value = -2050
$ws.Range("H61").Value = 6267.975
$ws.Range("I61").Value = 4518.2593
$ws.Range("J61").Value = 9902
$ws.Range("K61").Value = 4518.2593
$ws.Range("L61").Value = 9902
$ws.Range("M61").Value = -4306.2593
$ws.Range("N61").Value = -10326
$ws.Range("H132").Value = 4355.173
$ws.Range("I132").Value = 1789.6666
$ws.Range("J132").Value = 8811.053
$ws.Range("K132").Value = 5368.9998
$ws.Range("L132").Value = 26433.159
$ws.Range("M132").Value = -2838.9998
$ws.Range("N132").Value = -31493.159
$ws.Range("H136").Value = 6267.975
$ws.Range("I136").Value = 4518.2593
$ws.Range("J136").Value = 9902
$ws.Range("K136").Value = 13554.7779
$ws.Range("L136").Value = 29706
$ws.Range("M136").Value = -11004.7779
$ws.Range("N136").Value = -34806

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 309.1
$ws.Range("I22").Value = 230.33333
$ws.Range("J22").Value = 342.85715
$ws.Range("K22").Value = 230.33333
$ws.Range("L22").Value = 342.85715
$ws.Range("M22").Value = -57.33332999999999
$ws.Range("N22").Value = -688.85715
$ws.Range("H30").Value = 600
$ws.Range("I30").Value = 600
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 600
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -475
$ws.Range("N30").ClearContents()
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H36").Value = 3722.8333
$ws.Range("I36").Value = 884.25
$ws.Range("J36").Value = 9400
$ws.Range("K36").Value = 884.25
$ws.Range("L36").Value = 9400
$ws.Range("M36").Value = -350.25
$ws.Range("N36").Value = -10468
$ws.Range("H134").Value = 3443.1738
$ws.Range("I134").Value = 3791.3572
$ws.Range("J134").Value = 2901.5557
$ws.Range("K134").Value = 11374.0716
$ws.Range("L134").Value = 8704.667099999999
$ws.Range("M134").Value = -8839.071599999999
$ws.Range("N134").Value = -13774.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 120.2
$ws.Range("I22").Value = 125.25
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 125.25
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 224.75
$ws.Range("N22").Value = -800
$ws.Range("H36").Value = 13666.667
$ws.Range("I36").Value = 1000
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 1000
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -612
$ws.Range("N36").Value = -20776
$ws.Range("H40").Value = 13666.667
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 20000
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 20000
$ws.Range("M40").Value = -840
$ws.Range("N40").Value = -20320
$ws.Range("H132").Value = 2326.658
$ws.Range("I132").Value = 1796.8518
$ws.Range("J132").Value = 3627.0908
$ws.Range("K132").Value = 5390.555399999999
$ws.Range("L132").Value = 10881.2724
$ws.Range("M132").Value = -2860.555399999999
$ws.Range("N132").Value = -15941.2724
$ws.Range("H134").Value = 3499.9285
$ws.Range("I134").Value = 2057.85
$ws.Range("J134").Value = 4810.909
$ws.Range("K134").Value = 6173.549999999999
$ws.Range("L134").Value = 14432.727
$ws.Range("M134").Value = -3638.549999999999
$ws.Range("N134").Value = -19502.727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H113").Value = 3231.6365
$ws.Range("I113").Value = 3118.625
$ws.Range("J113").Value = 3533
$ws.Range("K113").Value = 3118.625
$ws.Range("L113").Value = 3533
$ws.Range("M113").Value = -948.625
$ws.Range("N113").Value = -7873

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 22700
$ws.Range("J94").Value = 22700
$ws.Range("L94").Value = 22700
$ws.Range("N94").Value = -24052
$ws.Range("H132").Value = 6377.4585
$ws.Range("I132").Value = 8772.346
$ws.Range("J132").Value = 3547.1365
$ws.Range("K132").Value = 26317.038
$ws.Range("L132").Value = 10641.4095
$ws.Range("M132").Value = -23787.038
$ws.Range("N132").Value = -15701.4095
$ws.Range("H136").Value = 3829.3394
$ws.Range("I136").Value = 2073.3713
$ws.Range("J136").Value = 6755.952
$ws.Range("K136").Value = 6220.113899999999
$ws.Range("L136").Value = 20267.856
$ws.Range("M136").Value = -3670.113899999999
$ws.Range("N136").Value = -25367.856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1000
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 373
$ws.Range("N96").ClearContents()
$ws.Range("H132").Value = 1976.7234
$ws.Range("I132").Value = 1000.6667
$ws.Range("J132").Value = 3294.4
$ws.Range("K132").Value = 3002.0001
$ws.Range("L132").Value = 9883.200000000001
$ws.Range("M132").Value = -472.0001000000002
$ws.Range("N132").Value = -14943.2
$ws.Range("H136").Value = 5215.532
$ws.Range("I136").Value = 3570.907
$ws.Range("J136").Value = 8937.579
$ws.Range("K136").Value = 10712.721
$ws.Range("L136").Value = 26812.737
$ws.Range("M136").Value = -8162.721000000001
$ws.Range("N136").Value = -31912.737
